$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-column formatting (style index 2, date number format) from the
# last existing data row (233) down across the new rows (234-238), matching
# column A's style in the diff.
$ws.Range("A233").Copy()
$ws.Range("A234:A238").PasteSpecial(-4122)

# New data rows appended to the report ("aggiornamento fino al 26/03"):
# date serials 44308-44312 => 2021-04-22 .. 2021-04-26, with their
# new-positives / 7-day-rolling-sum / rolling-sum-per-100k columns.
$data = @(
    @(234, 44308, 1, 5, 31.30870381966186),
    @(235, 44309, 4, 9, 56.35566687539136),
    @(236, 44310, 1, 8, 50.09392611145898),
    @(237, 44311, 3, 10, 62.61740763932373),
    @(238, 44312, 2, 12, 75.14088916718849)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
